$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.9792566666666667
$ws.Range("H2").Value = 2.93777
$ws.Range("I2").Value = 0.02840115057834171
$ws.Range("J2").Value = 0.02840115057834171
$ws.Range("M2").Value = 6.913788666666666
$ws.Range("N2").Value = 20.741366
$ws.Range("O2").Value = 0.4464851245108818
$ws.Range("P2").Value = 0.4464851245108818
$ws.Range("Q2").Value = 6.770373643757777
$ws.Range("R2").Value = 60.93336279382
$ws.Range("S2").Value = 0.0126806912522232
$ws.Range("T2").Value = 0.0126806912522232
# Row 3
$ws.Range("G3").Value = 0.9792566666666667
$ws.Range("H3").Value = 2.93777
$ws.Range("I3").Value = 0.02840115057834171
$ws.Range("J3").Value = 0.02840115057834171
$ws.Range("O3").Value = 0.002033181734278123
$ws.Range("P3").Value = 0.002033181734278123
$ws.Range("Q3").Value = 0.03083059047444444
$ws.Range("R3").Value = 0.27747531427
$ws.Range("S3").Value = 0.00005774470058836693
$ws.Range("T3").Value = 0.00005774470058836692
# Row 4
$ws.Range("G4").Value = 0.9792566666666667
$ws.Range("H4").Value = 2.93777
$ws.Range("I4").Value = 0.02840115057834171
$ws.Range("J4").Value = 0.02840115057834171
$ws.Range("M4").Value = 8.539652666666667
$ws.Range("N4").Value = 25.618958
$ws.Range("O4").Value = 0.55148169375484
$ws.Range("P4").Value = 0.55148169375484
$ws.Range("Q4").Value = 8.362511804851112
$ws.Range("R4").Value = 75.26260624366
$ws.Range("S4").Value = 0.01566271462553014
$ws.Range("T4").Value = 0.01566271462553014
# Row 5
$ws.Range("I5").Value = 0.06296007145894493
$ws.Range("J5").Value = 0.06296007145894492
$ws.Range("M5").Value = 6.913788666666666
$ws.Range("N5").Value = 20.741366
$ws.Range("O5").Value = 0.4464851245108818
$ws.Range("P5").Value = 0.4464851245108818
$ws.Range("Q5").Value = 15.00865985126
$ws.Range("R5").Value = 135.07793866134
$ws.Range("S5").Value = 0.02811073534456104
$ws.Range("T5").Value = 0.02811073534456104
# Row 6
$ws.Range("I6").Value = 0.06296007145894493
$ws.Range("J6").Value = 0.06296007145894492
$ws.Range("O6").Value = 0.002033181734278123
$ws.Range("P6").Value = 0.002033181734278123
$ws.Range("S6").Value = 0.0001280092672791722
$ws.Range("T6").Value = 0.0001280092672791722
# Row 7
$ws.Range("I7").Value = 0.06296007145894493
$ws.Range("J7").Value = 0.06296007145894492
$ws.Range("M7").Value = 8.539652666666667
$ws.Range("N7").Value = 25.618958
$ws.Range("O7").Value = 0.55148169375484
$ws.Range("P7").Value = 0.55148169375484
$ws.Range("Q7").Value = 18.53813419838
$ws.Range("R7").Value = 166.84320778542
$ws.Range("S7").Value = 0.03472132684710471
$ws.Range("T7").Value = 0.0347213268471047
# Row 8
$ws.Range("G8").Value = 18.980972
$ws.Range("H8").Value = 56.942916
$ws.Range("I8").Value = 0.550500662640664
$ws.Range("J8").Value = 0.550500662640664
$ws.Range("M8").Value = 6.913788666666666
$ws.Range("N8").Value = 20.741366
$ws.Range("O8").Value = 0.4464851245108818
$ws.Range("P8").Value = 0.4464851245108818
$ws.Range("Q8").Value = 131.2304290959173
$ws.Range("R8").Value = 1181.073861863256
$ws.Range("S8").Value = 0.2457903569024398
$ws.Range("T8").Value = 0.2457903569024398
# Row 9
$ws.Range("G9").Value = 18.980972
$ws.Range("H9").Value = 56.942916
$ws.Range("I9").Value = 0.550500662640664
$ws.Range("J9").Value = 0.550500662640664
$ws.Range("O9").Value = 0.002033181734278123
$ws.Range("P9").Value = 0.002033181734278123
$ws.Range("Q9").Value = 0.5975905954573333
$ws.Range("R9").Value = 5.378315359115999
$ws.Range("S9").Value = 0.001119267891989001
$ws.Range("T9").Value = 0.001119267891989001
# Row 10
$ws.Range("G10").Value = 18.980972
$ws.Range("H10").Value = 56.942916
$ws.Range("I10").Value = 0.550500662640664
$ws.Range("J10").Value = 0.550500662640664
$ws.Range("M10").Value = 8.539652666666667
$ws.Range("N10").Value = 25.618958
$ws.Range("O10").Value = 0.55148169375484
$ws.Range("P10").Value = 0.55148169375484
$ws.Range("Q10").Value = 162.0909081557253
$ws.Range("R10").Value = 1458.818173401528
$ws.Range("S10").Value = 0.3035910378462351
$ws.Range("T10").Value = 0.3035910378462351
# Row 11
$ws.Range("G11").Value = 0.9440396666666667
$ws.Range("H11").Value = 2.832119
$ws.Range("I11").Value = 0.02737976021771022
$ws.Range("J11").Value = 0.02737976021771022
$ws.Range("M11").Value = 6.913788666666666
$ws.Range("N11").Value = 20.741366
$ws.Range("O11").Value = 0.4464851245108818
$ws.Range("P11").Value = 0.4464851245108818
$ws.Range("Q11").Value = 6.526890748283778
$ws.Range("R11").Value = 58.742016734554
$ws.Range("S11").Value = 0.01222465564988244
$ws.Range("T11").Value = 0.01222465564988244
# Row 12
$ws.Range("G12").Value = 0.9440396666666667
$ws.Range("H12").Value = 2.832119
$ws.Range("I12").Value = 0.02737976021771022
$ws.Range("J12").Value = 0.02737976021771022
$ws.Range("O12").Value = 0.002033181734278123
$ws.Range("P12").Value = 0.002033181734278123
$ws.Range("Q12").Value = 0.02972183018544445
$ws.Range("R12").Value = 0.267496471669
$ws.Range("S12").Value = 0.00005566802836356324
$ws.Range("T12").Value = 0.00005566802836356323
# Row 13
$ws.Range("G13").Value = 0.9440396666666667
$ws.Range("H13").Value = 2.832119
$ws.Range("I13").Value = 0.02737976021771022
$ws.Range("J13").Value = 0.02737976021771022
$ws.Range("M13").Value = 8.539652666666667
$ws.Range("N13").Value = 25.618958
$ws.Range("O13").Value = 0.55148169375484
$ws.Range("P13").Value = 0.55148169375484
$ws.Range("Q13").Value = 8.061770856889112
$ws.Range("R13").Value = 72.55593771200201
$ws.Range("S13").Value = 0.01509943653946422
$ws.Range("T13").Value = 0.01509943653946422
# Row 14
$ws.Range("G14").Value = 11.40437333333333
$ws.Range("H14").Value = 34.21312
$ws.Range("I14").Value = 0.3307583551043392
$ws.Range("J14").Value = 0.3307583551043392
$ws.Range("M14").Value = 6.913788666666666
$ws.Range("N14").Value = 20.741366
$ws.Range("O14").Value = 0.4464851245108818
$ws.Range("P14").Value = 0.4464851245108818
$ws.Range("Q14").Value = 78.84742710243556
$ws.Range("R14").Value = 709.62684392192
$ws.Range("S14").Value = 0.1476786853617753
$ws.Range("T14").Value = 0.1476786853617753
# Row 15
$ws.Range("G15").Value = 11.40437333333333
$ws.Range("H15").Value = 34.21312
$ws.Range("I15").Value = 0.3307583551043392
$ws.Range("J15").Value = 0.3307583551043392
$ws.Range("O15").Value = 0.002033181734278123
$ws.Range("P15").Value = 0.002033181734278123
$ws.Range("Q15").Value = 0.3590514885688889
$ws.Range("R15").Value = 3.23146339712
$ws.Range("S15").Value = 0.0006724918460580198
$ws.Range("T15").Value = 0.0006724918460580196
# Row 16
$ws.Range("G16").Value = 11.40437333333333
$ws.Range("H16").Value = 34.21312
$ws.Range("I16").Value = 0.3307583551043392
$ws.Range("J16").Value = 0.3307583551043392
$ws.Range("M16").Value = 8.539652666666667
$ws.Range("N16").Value = 25.618958
$ws.Range("O16").Value = 0.55148169375484
$ws.Range("P16").Value = 0.55148169375484
$ws.Range("Q16").Value = 97.38938714766223
$ws.Range("R16").Value = 876.5044843289601
$ws.Range("S16").Value = 0.1824071778965058
$ws.Range("T16").Value = 0.1824071778965058
